# Revert "Powerpoint writer: consolidate text run nodes."
#
# The title shape ("A Table, with a caption") and the caption textbox
# ("Demonstration of simple table syntax, with alignment") each had their
# runs consolidated so that every run ended in a trailing space (e.g.
# "A ", "Table, ", "with ", ...). This restores the original, unconsolidated
# layout where each word and each inter-word space is its own <a:r> run.
#
# The host's TextRange.Text setter performs a minimal text diff against the
# existing runs and only touches what changed, so by first deleting just the
# space between two words (merging them into a single run) and then
# re-inserting that same space, the space is forced back out into its own,
# separate run - without disturbing any of the other, already-correct runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Split-Runs($tr, [string[]]$words) {
    $full = [string]::Join(" ", $words)
    # Make sure we start from the fully consolidated sentence.
    $tr.Text = $full
    for ($i = 0; $i -lt $words.Length - 1; $i++) {
        $squeezed = @()
        for ($j = 0; $j -lt $words.Length; $j++) {
            if ($j -eq $i + 1) {
                $squeezed[$squeezed.Length - 1] = $squeezed[$squeezed.Length - 1] + $words[$j]
            } else {
                $squeezed += $words[$j]
            }
        }
        $tr.Text = [string]::Join(" ", $squeezed)
        $tr.Text = $full
    }
}

$titleRange = $s.Shapes.Item("Title 1").TextFrame.TextRange
Split-Runs $titleRange @("A", "Table,", "with", "a", "caption")

$captionRange = $s.Shapes.Item("TextBox 3").TextFrame.TextRange
Split-Runs $captionRange @("Demonstration", "of", "simple", "table", "syntax,", "with", "alignment")
